$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57 (shifts rows 57..124 down to 58..125)
$ws.Rows.Item(57).Insert()

# Fill in the new row 57 with the new weekly record (same market/product
# metadata as every other row in this sheet; only the date/volume/price/
# origin columns differ row to row).
$ws.Range("A57").Value = 11
$ws.Range("B57").Value = "Vega Monumental Concepción"
$ws.Range("C57").Value = "Bíobío"
$ws.Range("D57").Value = 44792
$ws.Range("E57").Value = 8
$ws.Range("F57").Value = "Fruta"
$ws.Range("G57").Value = 100108
$ws.Range("H57").Value = "Tropicales y subtropicales"
$ws.Range("I57").Value = 100108002
$ws.Range("J57").Value = "Mango"
$ws.Range("K57").Value = "Sin especificar"
$ws.Range("L57").Value = "Primera"
$ws.Range("M57").Value = 100
$ws.Range("N57").Value = 9000
$ws.Range("O57").Value = 10000
$ws.Range("P57").Value = 9500
$ws.Range("Q57").Value = "$/bandeja 4 kilos"
$ws.Range("R57").Value = "Brasil"
$ws.Range("S57").Value = 2375
$ws.Range("T57").Value = 4
